$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.707.56'
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.804.90'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.34'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5923'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2775'
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06822'
$ws.Range("E9").Value = '  -3.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.34'
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07511'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.802.29'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.760'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6225'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.049.78'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009229'
$ws.Range("E16").Value = '  -6.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '75.67'
$ws.Range("E17").Value = '  -4.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.677.16'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("E19").Value = '  -6.35%  '
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '210.86'
$ws.Range("E21").Value = '  -6.71%  '
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.829'
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.93'
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.870'
$ws.Range("E26").Value = '  -1.77%  '
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.42'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("E29").Value = '  -3.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06162'
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.428'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.787'
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.743'
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.737'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.061'
$ws.Range("E35").Value = '  -5.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6428'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.498'
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.718'
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.549'
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01693'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.146.93'
$ws.Range("E41").Value = '  -5.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8845'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.94'
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.955.46'
$ws.Range("E45").Value = '  -1.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.59'
$ws.Range("E46").Value = '  -3.38%  '
$ws.Range("E47").Value = '  -2.25%  '
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.353'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05462'
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4475'
$ws.Range("E51").Value = '  -1.65%  '
